$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $ws.Range($cellRef).Value = "'" + $val
    $ws.Range($cellRef).Style = "Normal"
}

# Row 2
Set-TextValue "D2" "30.705.06"
$ws.Range("E2").Value = "  +2.03%  "

# Row 3
Set-TextValue "D3" "2.110.47"
$ws.Range("E3").Value = "  +10.27%  "

# Row 4
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
Set-TextValue "D5" "330.70"
$ws.Range("E5").Value = "  +3.21%  "

# Row 6
Set-TextValue "D6" "0.9997"
$ws.Range("E6").Value = "  -0.06%  "

# Row 7
Set-TextValue "D7" "0.5222"
$ws.Range("E7").Value = "  +3.06%  "

# Row 8
Set-TextValue "D8" "0.4410"
$ws.Range("E8").Value = "  +8.10%  "

# Row 9
Set-TextValue "D9" "0.09022"

# Row 10
Set-TextValue "D10" "46.95"
$ws.Range("E10").Value = "  +10.50%  "

# Row 11
Set-TextValue "D11" "1.179"
$ws.Range("E11").Value = "  +6.51%  "

# Row 12
Set-TextValue "D12" "25.16"
$ws.Range("E12").Value = "  +5.59%  "

# Row 13
Set-TextValue "D13" "2.110.59"
$ws.Range("E13").Value = "  +10.85%  "

# Row 14
Set-TextValue "D14" "6.793"
$ws.Range("E14").Value = "  +6.14%  "

# Row 15
Set-TextValue "D15" "7.746"
$ws.Range("E15").Value = "  +7.11%  "

# Row 16
Set-TextValue "D16" "97.96"
$ws.Range("E16").Value = "  +5.91%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D17" "0.00001140"
$ws.Range("E17").Value = "  +4.00%  "

# Row 18
$ws.Range("B18").Value = "BinanceUSD"
$ws.Range("C18").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D18" "1.000"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19
Set-TextValue "D19" "0.06625"
$ws.Range("E19").Value = "  +1.78%  "

# Row 20
Set-TextValue "D20" "19.27"
$ws.Range("E20").Value = "  +4.19%  "

# Row 21
$ws.Range("E21").Value = "  +8.03%  "

# Row 22
Set-TextValue "D22" "0.9996"

# Row 23
Set-TextValue "D23" "30.838.17"
$ws.Range("E23").Value = "  +2.44%  "

# Row 24
$ws.Range("E24").Value = "  +6.25%  "

# Row 25
Set-TextValue "D25" "2.356.94"
$ws.Range("E25").Value = "  +10.91%  "

# Row 26
Set-TextValue "D26" "2.261"
$ws.Range("E26").Value = "  +3.21%  "

# Row 27
Set-TextValue "D27" "22.99"
$ws.Range("E27").Value = "  +5.31%  "

# Row 28
Set-TextValue "D28" "2.544"
$ws.Range("E28").Value = "  +11.64%  "

# Row 29
Set-TextValue "D29" "163.27"
$ws.Range("E29").Value = "  +0.20%  "

# Row 30
Set-TextValue "D30" "133.99"
$ws.Range("E30").Value = "  +4.03%  "

# Row 31
Set-TextValue "D31" "1.191"

# Row 32
Set-TextValue "D32" "0.1069"
$ws.Range("E32").Value = "  +2.53%  "

# Row 33
Set-TextValue "D33" "6.239"
$ws.Range("E33").Value = "  +4.67%  "

# Row 34
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D34" "3.907"
$ws.Range("E34").Value = "  +2.85%  "

# Row 35
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D35" "1.531"
$ws.Range("E35").Value = "  +28.09%  "

# Row 36
Set-TextValue "D36" "0.02586"
$ws.Range("E36").Value = "  +5.22%  "

# Row 37
Set-TextValue "D37" "5.631"
$ws.Range("E37").Value = "  +5.07%  "

# Row 38
Set-TextValue "D38" "0.06750"
$ws.Range("E38").Value = "  +5.49%  "

# Row 39
$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D39" "12.82"
$ws.Range("E39").Value = "  +12.42%  "

# Row 40
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "9.530"
$ws.Range("E40").Value = "  +10.49%  "

# Row 41
Set-TextValue "D41" "0.2258"
$ws.Range("E41").Value = "  +4.84%  "

# Row 42
Set-TextValue "D42" "0.6799"
$ws.Range("E42").Value = "  +3.92%  "

# Row 43
Set-TextValue "D43" "1.255"
$ws.Range("E43").Value = "  +3.31%  "

# Row 44
Set-TextValue "D44" "14.35"
$ws.Range("E44").Value = "  +6.75%  "

# Row 45
$ws.Range("E45").Value = "  -0.01%  "

# Row 46
Set-TextValue "D46" "0.6339"
$ws.Range("E46").Value = "  +4.19%  "

# Row 47
Set-TextValue "D47" "2.254"
$ws.Range("E47").Value = "  +2.78%  "

# Row 48
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
Set-TextValue "D48" "1.290"
$ws.Range("E48").Value = "  +6.61%  "

# Row 49
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D49" "3.652"
$ws.Range("E49").Value = "  +0.77%  "

# Row 50
Set-TextValue "D50" "124.07"
$ws.Range("E50").Value = "  +1.60%  "

# Row 51
Set-TextValue "D51" "83.12"
$ws.Range("E51").Value = "  +5.15%  "
